$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.633.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.645.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.687.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.280.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("E36").Value = "  +6.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.786.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.19%  "
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0981"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
